$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions (style matches the existing header cells, e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows for new columns I (I0) and J (IF)
$data = @(
    @(2, 9, 9),
    @(3, 8, 8),
    @(4, 7, 9),
    @(5, 5, 6),
    @(6, 5, 5),
    @(7, 8, 8),
    @(8, 6, 7),
    @(9, 8, 8),
    @(10, 1, 6),
    @(11, 6, 6),
    @(12, 7, 7),
    @(13, 7, 7),
    @(14, 5, 5)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
